$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.973.28"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "1.655.18"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +2.58%  "
$ws.Range("D9").Value = "'0.0617"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").Value = "'20.16"
$ws.Range("E10").Value = "  +4.71%  "
$ws.Range("D11").Value = "'0.0879"
$ws.Range("E11").Value = "  +3.56%  "
$ws.Range("D12").Value = "1.887.83"
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").Value = "1.649.69"
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("E14").Value = "  +2.11%  "
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "'65.29"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").Value = "26.979.34"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "'237.13"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "0.0₃0737"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "'7.76"
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +3.80%  "
$ws.Range("D23").Value = "'9.29"
$ws.Range("E23").Value = "  +3.11%  "
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "'145.25"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'7.13"
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "'15.85"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").Value = "1.551.49"
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("E34").Value = "  +4.34%  "
$ws.Range("E35").Value = "  +10.38%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("E38").Value = "  +9.31%  "
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("E40").Value = "  +3.71%  "
$ws.Range("D42").Value = "'66.76"
$ws.Range("E42").Value = "  +9.45%  "
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").Value = "'0.966"
$ws.Range("E44").Value = "  +4.31%  "
$ws.Range("D45").Value = "1.796.45"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("D47").Value = "'90.04"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "'7.66"
$ws.Range("E51").Value = "  +2.88%  "
